# Correcting error in calculation of T* and adaptation analysis windows
# where windows were double the necessary length, due to miscalculation
# about frame rate (1f/s instead of 2f/s).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Halve/correct the existing T* (D column) window values that were
#     computed assuming 2 frames/sec instead of 1 frame/sec ---
$ws.Range("D4").Value  = 32.300000000000004
$ws.Range("D5").Value  = 31.5
$ws.Range("D9").Value  = 35
$ws.Range("D11").Value = 29.5
$ws.Range("D14").Value = 22.55
$ws.Range("D15").Value = 25.45
$ws.Range("D25").Value = 35.799999999999997
$ws.Range("D29").Value = 35.700000000000003
$ws.Range("D55").Value = 24.15
$ws.Range("D64").Value = 23.2
$ws.Range("D89").Value = 28
$ws.Range("D107").Value = 28.05

# --- Fill in previously-blank T*/Tmax (D/E) pairs now that the
#     analysis window calculation has been corrected ---
$ws.Range("D22").Value = 39.4
$ws.Range("E22").Value = 40.200000000000003

$ws.Range("D27").Value = 40.1
$ws.Range("E27").Value = 40

$ws.Range("D35").Value = 39.75
$ws.Range("E35").Value = 40.1

$ws.Range("D39").Value = 24.2
$ws.Range("E39").Value = 40

$ws.Range("D40").Value = 40
$ws.Range("E40").Value = 40.1

$ws.Range("D52").Value = 40
$ws.Range("E52").Value = 40

# --- Clear out the stale ASE_ctrl/XL115 Exp.UID entries (C74:C88) ---
#     These experiments' recording IDs were entered in error and are
#     removed entirely (leaving the cells blank, not zero).
$ws.Range("C74:C88").ClearContents()

# --- Update the saved selection/scroll position on Sheet1 ---
$ws.Activate()
$ws.Range("C2:E20").Select()
